$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Model")
$ws.Range("U9").Value = 22140
$wb.Save()
